# The "added transmission expansion to technology term" change flips the
# sign convention of the beta_ij curtailment-response coefficients: every
# populated cell in the beta_ij (SDES) table (rows 42-48) and the
# beta_ij (LDES) table (rows 51-57), columns D:J, on both the
# "wind curtailment" and "solar curtailment" sheets, is negated in place.
# Cells that were blank stay blank - only cells that already held a
# coefficient are rewritten, with the sign flipped.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("wind curtailment")
# row 42
$ws.Range("E42").Value = 0.03807920850308791
$ws.Range("F42").Value = 0.04926759243704416
$ws.Range("G42").Value = 0.1749825512014097
$ws.Range("H42").Value = 0.4657837775029747
$ws.Range("I42").Value = 0.3822968326306992
$ws.Range("J42").Value = 0.581644751777096
# row 43
$ws.Range("D43").Value = 0.01419223246491806
$ws.Range("E43").Value = 0.1410778478247
$ws.Range("F43").Value = 0.2251472196146845
$ws.Range("G43").Value = 0.3313660281446567
$ws.Range("H43").Value = 0.4351045304287916
$ws.Range("I43").Value = 0.5050387678763628
$ws.Range("J43").Value = 0.5553449424733008
# row 44
$ws.Range("D44").Value = 0.01624199507771802
$ws.Range("E44").Value = 0.1243593841431062
$ws.Range("F44").Value = 0.1782162170375382
$ws.Range("G44").Value = 0.228181530141783
$ws.Range("H44").Value = 0.2765849466592284
$ws.Range("I44").Value = 0.3144327993924718
$ws.Range("J44").Value = 0.3003866169489345
# row 45
$ws.Range("D45").Value = 0.01743064528986662
$ws.Range("E45").Value = 0.0898449861760492
$ws.Range("F45").Value = 0.106429764677167
$ws.Range("G45").Value = 0.1190856003175532
$ws.Range("H45").Value = 0.1143818786250513
$ws.Range("I45").Value = 0.1430240986471133
# row 46
$ws.Range("D46").Value = 0.01107935198956332
$ws.Range("E46").Value = 0.05180637406524109
$ws.Range("F46").Value = 0.04515241954845074
$ws.Range("G46").Value = 0.03699594123370264
$ws.Range("H46").Value = 0.007988524135573341
$ws.Range("I46").Value = 0.01591847493247362
# row 47
$ws.Range("D47").Value = 0.005953797437680362
$ws.Range("E47").Value = 0.01871502997681276
$ws.Range("F47").Value = 0.004549167439444669
$ws.Range("G47").Value = -0.03014924825547636
$ws.Range("H47").Value = -0.06856167593237956
# row 48
$ws.Range("D48").Value = 0.001004438729904129
$ws.Range("E48").Value = -0.02035137001495592
$ws.Range("F48").Value = -0.05788503803984875
# row 51
$ws.Range("E51").Value = 0.1510731638907951
$ws.Range("F51").Value = 0.2686526591317246
$ws.Range("G51").Value = 0.5130150794351762
$ws.Range("H51").Value = 0.6599230583560466
$ws.Range("I51").Value = 0.7761945502660924
$ws.Range("J51").Value = 1.084637948158285
# row 52
$ws.Range("D52").Value = 0.001707913159101493
$ws.Range("E52").Value = 0.2539799277840126
$ws.Range("F52").Value = 0.6461515859087472
$ws.Range("G52").Value = 0.8739608739055509
$ws.Range("H52").Value = 1.121633891275849
$ws.Range("I52").Value = 1.289035448916887
$ws.Range("J52").Value = 1.444521156155731
# row 53
$ws.Range("D53").Value = 0.004746608915359804
$ws.Range("E53").Value = 0.2683431549302445
$ws.Range("F53").Value = 0.5511170877225552
$ws.Range("G53").Value = 0.8510947657467732
$ws.Range("H53").Value = 1.077936295793381
$ws.Range("I53").Value = 1.289670401863402
$ws.Range("J53").Value = 1.119346987099373
# row 54
$ws.Range("D54").Value = 0.04115880605720874
$ws.Range("E54").Value = 0.2097646454874569
$ws.Range("F54").Value = 0.5006630936993596
$ws.Range("G54").Value = 0.8274623592207895
$ws.Range("H54").Value = 1.061814793797657
$ws.Range("I54").Value = 1.240909828752011
# row 55
$ws.Range("D55").Value = 0.008558362859703005
$ws.Range("E55").Value = 0.1880232783644833
$ws.Range("F55").Value = 0.5170951321489407
$ws.Range("G55").Value = 0.7605002217647031
$ws.Range("H55").Value = 1.068166512772024
$ws.Range("I55").Value = 1.094311838397483
# row 56
$ws.Range("D56").Value = 0.01243107319904072
$ws.Range("E56").Value = 0.2523794122609848
$ws.Range("F56").Value = 0.4962889934261825
$ws.Range("G56").Value = 0.7651035169817634
$ws.Range("H56").Value = 1.057533487248302
# row 57
$ws.Range("D57").Value = -0.009400868951178839
$ws.Range("E57").Value = 0.2687730762096723
$ws.Range("F57").Value = 0.5075483500591331

$ws = $wb.Worksheets.Item("solar curtailment")
# row 42
$ws.Range("E42").Value = 0.0443934115864938
$ws.Range("F42").Value = 0.001336947024145993
$ws.Range("G42").Value = 0.03779790995251268
$ws.Range("H42").Value = 0.2192311907808228
$ws.Range("I42").Value = 0.2142824060512421
$ws.Range("J42").Value = 0.3109339512502715
# row 43
$ws.Range("D43").Value = -0.000007170673537677454
$ws.Range("E43").Value = 0.04072135331409409
$ws.Range("F43").Value = 0.07424914552703732
$ws.Range("G43").Value = 0.06654355573083219
$ws.Range("H43").Value = 0.1122825620799758
$ws.Range("I43").Value = 0.1896266843091846
$ws.Range("J43").Value = 0.3501382794911784
# row 44
$ws.Range("D44").Value = 0.00004017680576350262
$ws.Range("E44").Value = 0.02635004268088499
$ws.Range("F44").Value = 0.08416579132352937
$ws.Range("G44").Value = 0.1311721569005507
$ws.Range("H44").Value = 0.1290375611687204
$ws.Range("I44").Value = 0.1389185734277663
$ws.Range("J44").Value = 0.2176971956367648
# row 45
$ws.Range("D45").Value = 0.000073482761622889
$ws.Range("E45").Value = 0.0486254041774131
$ws.Range("F45").Value = 0.09222245664265515
$ws.Range("G45").Value = 0.1107917325037082
$ws.Range("H45").Value = 0.1403886811658372
$ws.Range("I45").Value = 0.07961127639709582
# row 46
$ws.Range("D46").Value = 0.0001800821543901399
$ws.Range("E46").Value = 0.04700598487038497
$ws.Range("F46").Value = 0.1134377019231578
$ws.Range("G46").Value = 0.1213465856756812
$ws.Range("H46").Value = 0.05223099867876894
$ws.Range("I46").Value = -0.2204365736448612
# row 47
$ws.Range("D47").Value = 0.00005221878187348652
$ws.Range("E47").Value = 0.05311464896388466
$ws.Range("F47").Value = 0.07875486248291225
$ws.Range("G47").Value = -0.002089449809821217
$ws.Range("H47").Value = -0.2262773911747978
# row 48
$ws.Range("D48").Value = 0.0001483887263475372
$ws.Range("E48").Value = 0.03927707482323817
$ws.Range("F48").Value = -0.06147543694844773
# row 51
$ws.Range("E51").Value = 0.2908328138323258
$ws.Range("F51").Value = 0.5096988584552465
$ws.Range("G51").Value = 0.7191890220212594
$ws.Range("H51").Value = 0.8732452277949271
$ws.Range("I51").Value = 0.9566724135629523
$ws.Range("J51").Value = 1.032973940284662
# row 52
$ws.Range("D52").Value = 0.0005140527130918539
$ws.Range("E52").Value = 0.320053997515194
$ws.Range("F52").Value = 0.5696852569605324
$ws.Range("G52").Value = 0.7515295361662166
$ws.Range("H52").Value = 0.8811848077588235
$ws.Range("I52").Value = 0.9694244255569329
$ws.Range("J52").Value = 1.061241830821753
# row 53
$ws.Range("D53").Value = 0.0005987242167493896
$ws.Range("E53").Value = 0.3099320992833279
$ws.Range("F53").Value = 0.5692310729836303
$ws.Range("G53").Value = 0.7675684915317852
$ws.Range("H53").Value = 0.9039136213075745
$ws.Range("I53").Value = 0.9952994996996426
$ws.Range("J53").Value = 1.100508264115412
# row 54
$ws.Range("D54").Value = 0.001236127776125065
$ws.Range("E54").Value = 0.2991250808544557
$ws.Range("F54").Value = 0.573643424832342
$ws.Range("G54").Value = 0.7838988623476733
$ws.Range("H54").Value = 0.9295001512731182
$ws.Range("I54").Value = 1.034899869980614
# row 55
$ws.Range("D55").Value = 0.001239739247778578
$ws.Range("E55").Value = 0.2911804032713267
$ws.Range("F55").Value = 0.5816092507080989
$ws.Range("G55").Value = 0.8071614649744907
$ws.Range("H55").Value = 0.9596811512954909
$ws.Range("I55").Value = 1.066464881058169
# row 56
$ws.Range("D56").Value = 0.002206986528664534
$ws.Range("E56").Value = 0.2773199247978722
$ws.Range("F56").Value = 0.5975804945647406
$ws.Range("G56").Value = 0.828242120952577
$ws.Range("H56").Value = 0.9711025010416207
# row 57
$ws.Range("D57").Value = 0.004687280317108817
$ws.Range("E57").Value = 0.2574871344819252
$ws.Range("F57").Value = 0.6079428170164718

